$wb = $excel.ActiveWorkbook

# --- Overview sheet: mark the 09458abb... file row as handed back ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: record the handback target/file/datetime for 09458abb... ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("I2").Value = "09458abb-652e-48e8-8e5d-9581ec1c5232.md"
$zhcn.Range("I2").Style = "Hyperlink"
$zhcn.Range("J2").Value = "09458abb-652e-48e8-8e5d-9581ec1c5232.31d4f97bc0e5eeeb11e0dc9f635ddbb294c944ef.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-18 10:44:11"

# --- de-de sheet: record the handback target/file/datetime for 09458abb... ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("I2").Value = "09458abb-652e-48e8-8e5d-9581ec1c5232.md"
$dede.Range("I2").Style = "Hyperlink"
$dede.Range("J2").Value = "09458abb-652e-48e8-8e5d-9581ec1c5232.31d4f97bc0e5eeeb11e0dc9f635ddbb294c944ef.de-de.xlf"
$dede.Range("K2").Value = "2016-08-18 10:44:19"
